$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append the new e-mail log row (row 33) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A33").Value = "Offerte voor 500 stuks"
$logs.Range("B33").Value = "mailmind.test@zohomail.eu"
$logs.Range("C33").Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$logs.Range("D33").Value = "Offerte / Prijsaanvraag"
$logs.Range("E33").Value = "Geachte klant,`nHartelijk dank voor uw interesse in product X. Om u een nauwkeurige offerte te kunnen sturen, hebben we wat meer informatie nodig. Kunt u ons laten weten of u specifieke eisen heeft met betrekking tot het product of de levering? Ook is het handig als u de gewenste leverdatum kunt vermelden.`nZodra we deze details van u ontvangen, zullen we zo spoedig mogelijk een offerte voor 500 stuks van product X opstellen.`nMet vriendelijke groet,`n[Naam van het bedrijf] E-mailassistent"
$logs.Range("F33").Value = "2025-06-22 19:03:12"
$logs.Range("G33").Value = "Ja"

# Row 33 contains a multi-line value (column E); avoid leaving a stray
# custom row-height behind from the implicit wrap auto-sizing.
$logs.Rows.Item(33).AutoFit()

# Conditional formatting ranges grow to cover the new row.
$fcD = $logs.Range("D2:D32").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($logs.Range("D2:D33"))

$fcG = $logs.Range("G2:G32").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($logs.Range("G2:G33"))

# --- Sheet "Dashboard": category counts changed (row 7 / row 9 swap + increment) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Offerte / Prijsaanvraag"
$dash.Range("B7").Value = 3

$dash.Range("A9").Value = "Sollicitatie / Vacature"
$dash.Range("B9").Value = 2
